$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

# Insert a new row above row 69 (pushes existing rows 69-72 down to 70-73).
$ws.Rows.Item(69).Insert()

# Copy the formatting of the row above (row 68, "Hallinta"/"Integraatioloki")
# onto the freshly inserted row, matching Excel's own behaviour of carrying
# the look of the preceding row into a newly inserted one.
$ws.Range("A68:X68").Copy($ws.Range("A69:X69"))

# Fill in the new row's content: Hallinta / Yhteydenpito / (blank) / R*,W*
$ws.Range("A69").Value = "Hallinta"
$ws.Range("B69").Value = "Yhteydenpito"
$ws.Range("D69").Value = "R*,W*"
